# Update Crlf1-Cntfr.xlsx with new TPM-derived NATMI output.
# - The "Neutrophils" cluster is renamed to "Resolving-Mac".
# - A new cluster "Inflammatory-Mac" is introduced, taking over the
#   "target cluster" slot previously held by "MuSCs" in rows 4/8/12/16,
#   while "MuSCs" itself moves into the slot previously held by
#   "Neutrophils" in rows 5/9/13/17.
# - All numeric LR-pair statistics (columns G-T) are refreshed with the
#   newly computed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Neutrophils" cluster label to "Resolving-Mac" everywhere it appears
# (column A for rows 14-17, the sending-cluster column for that cluster)
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("A17").Value = "Resolving-Mac"

# Row 2
$ws.Range("G2").Value = 0.4287225
$ws.Range("H2").Value = 0.857445
$ws.Range("I2").Value = 0.04766385615009357
$ws.Range("J2").Value = 0.04112453143803989
$ws.Range("M2").Value = 0.08210149999999999
$ws.Range("N2").Value = 0.164203
$ws.Range("O2").Value = 0.01959206455542894
$ws.Range("P2").Value = 0.0143098977453608
$ws.Range("Q2").Value = 0.03519876033375
$ws.Range("R2").Value = 0.140795041335
$ws.Range("S2").Value = 0.000933833346653312
$ws.Range("T2").Value = 0.0005884878397042266

# Row 3
$ws.Range("G3").Value = 0.4287225
$ws.Range("H3").Value = 0.857445
$ws.Range("I3").Value = 0.04766385615009357
$ws.Range("J3").Value = 0.04112453143803989
$ws.Range("O3").Value = 0.7352800102707316
$ws.Range("P3").Value = 0.8055645487039627
$ws.Range("Q3").Value = 1.3209911996
$ws.Range("R3").Value = 7.925947197600001
$ws.Range("S3").Value = 0.03504628063958348
$ws.Range("T3").Value = 0.03312846460854653

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.4287225
$ws.Range("H4").Value = 0.857445
$ws.Range("I4").Value = 0.04766385615009357
$ws.Range("J4").Value = 0.04112453143803989
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.012461
$ws.Range("N4").Value = 0.037383
$ws.Range("O4").Value = 0.002973596297573126
$ws.Range("P4").Value = 0.003257838817895063
$ws.Range("Q4").Value = 0.0053423110725
$ws.Range("R4").Value = 0.032053866435
$ws.Range("S4").Value = 0.0001417330661759763
$ws.Range("T4").Value = 0.0001339770948865922

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.4287225
$ws.Range("H5").Value = 0.857445
$ws.Range("I5").Value = 0.04766385615009357
$ws.Range("J5").Value = 0.04112453143803989
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.0147595
$ws.Range("N5").Value = 2.029519
$ws.Range("O5").Value = 0.2421543288762665
$ws.Range("P5").Value = 0.1768677147327815
$ws.Range("Q5").Value = 0.4350502297387499
$ws.Range("R5").Value = 1.740200918955
$ws.Range("S5").Value = 0.01154200909768081
$ws.Range("T5").Value = 0.007273601894902543

# Row 6
$ws.Range("I6").Value = 0.3029995857564555
$ws.Range("J6").Value = 0.3921435547802385
$ws.Range("M6").Value = 0.08210149999999999
$ws.Range("N6").Value = 0.164203
$ws.Range("O6").Value = 0.01959206455542894
$ws.Range("P6").Value = 0.0143098977453608
$ws.Range("Q6").Value = 0.2237588533895
$ws.Range("R6").Value = 1.342553120337
$ws.Range("S6").Value = 0.005936387444408702
$ws.Range("T6").Value = 0.005611534170407506

# Row 7
$ws.Range("I7").Value = 0.3029995857564555
$ws.Range("J7").Value = 0.3921435547802385
$ws.Range("O7").Value = 0.7352800102707316
$ws.Range("P7").Value = 0.8055645487039627
$ws.Range("S7").Value = 0.222789538527034
$ws.Range("T7").Value = 0.3158969457337105

# Row 8
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.3029995857564555
$ws.Range("J8").Value = 0.3921435547802385
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.012461
$ws.Range("N8").Value = 0.037383
$ws.Range("O8").Value = 0.002973596297573126
$ws.Range("P8").Value = 0.003257838817895063
$ws.Range("Q8").Value = 0.03396112217300001
$ws.Range("R8").Value = 0.3056500995570001
$ws.Range("S8").Value = 0.0009009984463715869
$ws.Range("T8").Value = 0.00127754049495042

# Row 9
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.3029995857564555
$ws.Range("J9").Value = 0.3921435547802385
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.0147595
$ws.Range("N9").Value = 2.029519
$ws.Range("O9").Value = 0.2421543288762665
$ws.Range("P9").Value = 0.1768677147327815
$ws.Range("Q9").Value = 2.7656184379835
$ws.Range("R9").Value = 16.593710627901
$ws.Range("S9").Value = 0.07337266133864122
$ws.Range("T9").Value = 0.06935753438117009

# Row 10
$ws.Range("G10").Value = 5.7054395
$ws.Range("H10").Value = 11.410879
$ws.Range("I10").Value = 0.6343106498983883
$ws.Range("J10").Value = 0.5472853094614456
$ws.Range("M10").Value = 0.08210149999999999
$ws.Range("N10").Value = 0.164203
$ws.Range("O10").Value = 0.01959206455542894
$ws.Range("P10").Value = 0.0143098977453608
$ws.Range("Q10").Value = 0.4684251411092499
$ws.Range("R10").Value = 1.873700564437
$ws.Range("S10").Value = 0.01242745520100531
$ws.Range("T10").Value = 0.00783159681593143

# Row 11
$ws.Range("G11").Value = 5.7054395
$ws.Range("H11").Value = 11.410879
$ws.Range("I11").Value = 0.6343106498983883
$ws.Range("J11").Value = 0.5472853094614456
$ws.Range("O11").Value = 0.7352800102707316
$ws.Range("P11").Value = 0.8055645487039627
$ws.Range("Q11").Value = 17.57975233245334
$ws.Range("R11").Value = 105.47851399472
$ws.Range("S11").Value = 0.4663959411721214
$ws.Range("T11").Value = 0.4408736433286179

# Row 12
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 5.7054395
$ws.Range("H12").Value = 11.410879
$ws.Range("I12").Value = 0.6343106498983883
$ws.Range("J12").Value = 0.5472853094614456
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.012461
$ws.Range("N12").Value = 0.037383
$ws.Range("O12").Value = 0.002973596297573126
$ws.Range("P12").Value = 0.003257838817895063
$ws.Range("Q12").Value = 0.07109548160949999
$ws.Range("R12").Value = 0.426572889657
$ws.Range("S12").Value = 0.001886183800049051
$ws.Range("T12").Value = 0.00178296732562721

# Row 13
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 5.7054395
$ws.Range("H13").Value = 11.410879
$ws.Range("I13").Value = 0.6343106498983883
$ws.Range("J13").Value = 0.5472853094614456
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.0147595
$ws.Range("N13").Value = 2.029519
$ws.Range("O13").Value = 0.2421543288762665
$ws.Range("P13").Value = 0.1768677147327815
$ws.Range("Q13").Value = 5.789648934300248
$ws.Range("R13").Value = 23.15859573720099
$ws.Range("S13").Value = 0.1536010697252126
$ws.Range("T13").Value = 0.09679710199126898

# Row 14
$ws.Range("G14").Value = 0.1351536666666667
$ws.Range("H14").Value = 0.405461
$ws.Range("I14").Value = 0.01502590819506253
$ws.Range("J14").Value = 0.01944660432027604
$ws.Range("M14").Value = 0.08210149999999999
$ws.Range("N14").Value = 0.164203
$ws.Range("O14").Value = 0.01959206455542894
$ws.Range("P14").Value = 0.0143098977453608
$ws.Range("Q14").Value = 0.01109631876383333
$ws.Range("R14").Value = 0.066577912583
$ws.Range("S14").Value = 0.0002943885633616138
$ws.Range("T14").Value = 0.0002782789193176418

# Row 15
$ws.Range("G15").Value = 0.1351536666666667
$ws.Range("H15").Value = 0.405461
$ws.Range("I15").Value = 0.01502590819506253
$ws.Range("J15").Value = 0.01944660432027604
$ws.Range("O15").Value = 0.7352800102707316
$ws.Range("P15").Value = 0.8055645487039627
$ws.Range("Q15").Value = 0.4164390818311112
$ws.Range("R15").Value = 3.747951736480001
$ws.Range("S15").Value = 0.01104824993199265
$ws.Range("T15").Value = 0.0156654950330877

# Row 16
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 0.1351536666666667
$ws.Range("H16").Value = 0.405461
$ws.Range("I16").Value = 0.01502590819506253
$ws.Range("J16").Value = 0.01944660432027604
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.012461
$ws.Range("N16").Value = 0.037383
$ws.Range("O16").Value = 0.002973596297573126
$ws.Range("P16").Value = 0.003257838817895063
$ws.Range("Q16").Value = 0.001684149840333333
$ws.Range("R16").Value = 0.015157348563
$ws.Range("S16").Value = 0.00004468098497651164
$ws.Range("T16").Value = 0.00006335390243084113

# Row 17
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 0.1351536666666667
$ws.Range("H17").Value = 0.405461
$ws.Range("I17").Value = 0.01502590819506253
$ws.Range("J17").Value = 0.01944660432027604
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.0147595
$ws.Range("N17").Value = 2.029519
$ws.Range("O17").Value = 0.2421543288762665
$ws.Range("P17").Value = 0.1768677147327815
$ws.Range("Q17").Value = 0.1371484672098333
$ws.Range("R17").Value = 0.8228908032589999
$ws.Range("S17").Value = 0.00363858871473176
$ws.Range("T17").Value = 0.003439476465439859
